$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.84'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.23%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '29.75'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.40%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.284'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2.13%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05748'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.85%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.651'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.88%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.227'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '5.95%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8581'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.04%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-2.17%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1383'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.29%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07087'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.02%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03234'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '12.75%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09343'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.45%'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.33%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0005955'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.59%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005969'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.96%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.516'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.193'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-3.08%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03340'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.44%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1304'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.31%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.479'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '19.81%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.10%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.12%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004173'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-17.96%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.94%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-25.32%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03756'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.18%'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.11%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.002404'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '14.42%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-48.39%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009210'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-1.97%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005280'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '3.35%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000749'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.11%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.08981'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '26.46%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002191'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-19.11%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002098'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.11%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001998'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.11%'
